# Apply updated betting odds values to Sheet1 as per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7 updates
$ws.Range("M7").Value = 1.13
$ws.Range("N7").Value = 6

# Row 8 updates
$ws.Range("G8").Value = 1.44
$ws.Range("I8").Value = 7
$ws.Range("J8").Value = 2
$ws.Range("AH8").Value = 34
$ws.Range("AI8").Value = 21
$ws.Range("AN8").Value = 3.4
$ws.Range("AU8").Value = 9
$ws.Range("AW8").Value = 8
